# Foi excluído algumas funções não utilizadas
# Row 2 values were stored as text ("5000.0", "5000.0", "20"); convert them
# to real numbers. Add a new row 3 with the "old" text-style values
# ("4000.0", "4000.0", "50") that used to live in row 2, preserved as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: now numeric values
$ws.Range("A2").Value = 5000
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 20

# Row 3 (new): keep these as literal text, just like the old row 2 used to be.
# A leading apostrophe forces Excel to store a numeric-looking entry as text
# instead of silently converting it to a number.
$ws.Range("A3").Value = "'4000.0"
$ws.Range("B3").Value = "'4000.0"
$ws.Range("C3").Value = "'50"
